# EPBDS-13258 Property file is not read from dependent project
#
# The "Rules" sheet contains a small test table (rows 25-32, columns D:F)
# that shows the text resolved by the `msg(...)` user-defined function for
# three message sources: the child project's own MessageBundle (column D),
# a MessageBundle shipped inside a jar dependency (column E), and the
# parent/Main project resolving through its dependency (column F).
#
# Before the fix, message keys coming from a dependent/child project's
# property file were not resolved at all, so the raw lookup keys
# ("say.hello", "say.hello.1", "jar.say.hello", "jar.say.hello.1") were
# shown instead of the localized text. After the fix the real, resolved
# (and in the French-locale rows, translated) text is displayed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25 ("_res_.$Step1"): resolved key "say.hello" -> localized greeting.
$ws.Range("D25").Value = "Hello, from Project!"
$ws.Range("E25").Value = "Hello, from Project!"
$ws.Range("F25").Value = "Hello, from Project!"

# Row 29 ("_res_.$Step5"): resolved key "jar.say.hello" -> localized greeting.
$ws.Range("D29").Value = "Hello, from Jar!"
$ws.Range("E29").Value = "Hello, from Jar!"
$ws.Range("F29").Value = "Hello, from Jar!"

# Row 28 ("_res_.$Step4", French locale): jar message with parameter is now
# properly translated for the jar/MessageBundle column.
$ws.Range("E28").Value = "Bonjour, Parameter!"

# Row 27 ("_res_.$Step3", French locale): jar message without parameter is
# now properly translated for the jar/MessageBundle column.
$ws.Range("E27").Value = "Bonjour, from MessageBundle!"

# Row 26 ("_res_.$Step2"): resolved key "say.hello.1" -> localized greeting
# with parameter (same resolved text already used for other "with
# parameter" rows).
$ws.Range("D26").Value = "Hello, Parameter!"
$ws.Range("E26").Value = "Hello, Parameter!"
$ws.Range("F26").Value = "Hello, Parameter!"

# Row 30 ("_res_.$Step6"): resolved key "jar.say.hello.1" -> localized
# greeting with parameter.
$ws.Range("D30").Value = "Hello, Parameter!"
$ws.Range("E30").Value = "Hello, Parameter!"
$ws.Range("F30").Value = "Hello, Parameter!"

# Column D is now wide enough to comfortably display the resolved text
# (it used to share the narrow default column width).
$ws.Columns("D").ColumnWidth = 33

# Leave the selection on the cell the author was last reviewing.
$ws.Range("E17").Select()
